$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: set directly ---
$ws.Range('D2').Value = '65.411.74'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '2.950.52'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('E6').Value = '  +3.72%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').Value = '2.943.00'
$ws.Range('E10').Value = '  -4.88%  '
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('E13').Value = '  +1.78%  '
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '65.442.03'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '3.444.62'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '2.954.86'
$ws.Range('E19').Value = '  -1.81%  '
$ws.Range('E20').Value = '  +9.34%  '
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('E26').Value = '  -3.25%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E28').Value = '  -6.54%  '
$ws.Range('E29').Value = '  +2.57%  '
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('E37').Value = '  -0.86%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -7.55%  '
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E42').Value = '  -2.85%  '
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').Value = '2.712.32'
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('E50').Value = '  +5.54%  '
$ws.Range('E51').Value = '  +0.44%  '

# --- Numeric-looking Price values that must stay as literal text ---
# (Excel would otherwise auto-convert "1.00" -> 1, "0.690" -> 0.69, etc.)
# Temporarily force a Text number format, assign, then restore the default
# "Normal" style so no stray formatting is left behind.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '446.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.690'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.977'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '44.26'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.120'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.84'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.298'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '385.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0352'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.10'
$ws.Range('D48').Style = 'Normal'
